# Further WebAppData fixes for web tool scenario
#
# This reopens the "Boolean Disable Carbon Tax Border Adjustment" control-lever
# workbook, fixes a typo in the lever's description label, and leaves the
# BDCTBA sheet as the active/selected sheet (as it was when the file was last
# saved by the author), with cell C2 selected.

$wb = $excel.ActiveWorkbook
$wsAbout  = $wb.Worksheets.Item("About")
$wsBDCTBA = $wb.Worksheets.Item("BDCTBA")

# Fix the typo "Diable" -> "Disable" in the lever-name label on the BDCTBA sheet.
$wsBDCTBA.Range("B1").Value = "Disable Carbon Tax Border Adjustment"

# Make BDCTBA the active sheet, with C2 selected, matching the saved view state
# (the About sheet keeps its prior selection of A10).
$wsBDCTBA.Activate()
$wsBDCTBA.Range("C2").Select()

# Match the saved "do not recalculate before saving" workbook option.
$excel.CalculateBeforeSave = $false
